$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update D4 value (tiny precision change)
$ws.Range("D4").Value = 45689.69669267361

# Update D5 value
$ws.Range("D5").Value = 45689.74915857639

# Add new row 6: CUSTOM
$ws.Range("A6").Value = "CUSTOM"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 45689.73709944444
$ws.Range("D6").Value = 45689.74339052083

# Add new row 7: MAT151
$ws.Range("A7").Value = "MAT151"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = 45689.73989837963
$ws.Range("D7").Value = 45689.73992152778

# Add new row 8: MILLONIARA
$ws.Range("A8").Value = "MILLONIARA"
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = 45689.75070482639
$ws.Range("D8").Value = 45689.75072797792

# Copy style from existing data rows (C5:D5) down to new rows so that
# number formatting (style index 2, date/time) is preserved
$ws.Range("C5:D5").Copy()
$ws.Range("C6:D8").PasteSpecial(-4122)  # xlPasteFormats
